$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.563.82"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.560.89"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'210.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'0.515"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.41%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "'24.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.44%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").Value = "'0.0587"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "'0.0898"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "1.768.73"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "1.558.94"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "28.600.52"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "'61.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").Value = "'227.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").Value = "'7.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("D20").Value = "0.0₃0685"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").Value = "'0.997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").Value = "'3.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "'9.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.27%  "
$ws.Range("D24").Value = "'2.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.29%  "
$ws.Range("D25").Value = "'152.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").Value = "  +3.12%  "
$ws.Range("D27").Value = "'14.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").Value = "'6.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").Value = "'0.0458"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").Value = "'3.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "1.403.63"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "'3.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("D36").Value = "'1.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("D37").Value = "'2.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("D39").Value = "'0.0161"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'1.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Value = "'0.515"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").Value = "'0.770"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").Value = "'0.0461"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("D45").Value = "'63.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.00%  "
$ws.Range("D46").Value = "'5.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("D47").Value = "1.695.35"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").Value = "'0.832"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.15%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'84.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0103"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'41.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.45%  "
